$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)
$nl = [char]11

$cell = $tbl.Cell(1, 1)
$cell.Range.Text = "21 x 77" + $nl + "  7    7" + $nl + "  ----" + $nl + "2|    |" + $nl + "1|    |"

$cell = $tbl.Cell(1, 2)
$cell.Range.Text = "81 x 60" + $nl + "  6    0" + $nl + "  ----" + $nl + "8|    |" + $nl + "1|    |"

$cell = $tbl.Cell(1, 3)
$cell.Range.Text = "64 x 58" + $nl + "  5    8" + $nl + "  ----" + $nl + "6|    |" + $nl + "4|    |"

$cell = $tbl.Cell(2, 1)
$cell.Range.Text = "22 x 81" + $nl + "  8    1" + $nl + "  ----" + $nl + "2|    |" + $nl + "2|    |"

$cell = $tbl.Cell(2, 2)
$cell.Range.Text = "73 x 61" + $nl + "  6    1" + $nl + "  ----" + $nl + "7|    |" + $nl + "3|    |"

$cell = $tbl.Cell(2, 3)
$cell.Range.Text = "86 x 30" + $nl + "  3    0" + $nl + "  ----" + $nl + "8|    |" + $nl + "6|    |"

$cell = $tbl.Cell(3, 1)
$cell.Range.Text = "37 x 39" + $nl + "  3    9" + $nl + "  ----" + $nl + "3|    |" + $nl + "7|    |"

$cell = $tbl.Cell(3, 2)
$cell.Range.Text = "74 x 37" + $nl + "  3    7" + $nl + "  ----" + $nl + "7|    |" + $nl + "4|    |"

$cell = $tbl.Cell(3, 3)
$cell.Range.Text = "63 x 52" + $nl + "  5    2" + $nl + "  ----" + $nl + "6|    |" + $nl + "3|    |"

$cell = $tbl.Cell(4, 1)
$cell.Range.Text = "16 x 97" + $nl + "  9    7" + $nl + "  ----" + $nl + "1|    |" + $nl + "6|    |"

$cell = $tbl.Cell(4, 2)
$cell.Range.Text = "27 x 39" + $nl + "  3    9" + $nl + "  ----" + $nl + "2|    |" + $nl + "7|    |"

$cell = $tbl.Cell(4, 3)
$cell.Range.Text = "87 x 98" + $nl + "  9    8" + $nl + "  ----" + $nl + "8|    |" + $nl + "7|    |"

$cell = $tbl.Cell(5, 1)
$cell.Range.Text = "89 x 95" + $nl + "  9    5" + $nl + "  ----" + $nl + "8|    |" + $nl + "9|    |"

$cell = $tbl.Cell(5, 2)
$cell.Range.Text = "69 x 78" + $nl + "  7    8" + $nl + "  ----" + $nl + "6|    |" + $nl + "9|    |"

$cell = $tbl.Cell(5, 3)
$cell.Range.Text = "96 x 12" + $nl + "  1    2" + $nl + "  ----" + $nl + "9|    |" + $nl + "6|    |"
